$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2638.6667
$ws.Range("J17").Value = 2638.6667
$ws.Range("L17").Value = 7916.000100000001
$ws.Range("N17").Value = -8252.000100000001

$ws.Range("H41").Value = 729.44116
$ws.Range("I41").Value = 648.16
$ws.Range("J41").Value = 955.2222
$ws.Range("K41").Value = 648.16
$ws.Range("L41").Value = 955.2222
$ws.Range("M41").Value = -208.16
$ws.Range("N41").Value = -1835.2222

$ws.Range("H58").Value = 2478
$ws.Range("J58").Value = 3704
$ws.Range("L58").Value = 11112
$ws.Range("N58").Value = -11412

$ws.Range("H132").Value = 1055.3903
$ws.Range("I132").Value = 976.9394
$ws.Range("J132").Value = 1379
$ws.Range("K132").Value = 2930.8182
$ws.Range("L132").Value = 4137
$ws.Range("M132").Value = -400.8181999999997
$ws.Range("N132").Value = -9197

$ws.Range("H135").Value = 17247520
$ws.Range("I135").Value = 23811244
$ws.Range("K135").Value = 214301196
$ws.Range("M135").Value = -214298661

$ws.Range("H138").Value = 2804.3015
$ws.Range("J138").Value = 3351.1086
$ws.Range("L138").Value = 10053.3258
$ws.Range("N138").Value = -20333.3258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14096696
$ws.Range("I32").Value = 16960142
$ws.Range("J32").Value = 18080.916
$ws.Range("K32").Value = 16960142
$ws.Range("L32").Value = 18080.916
$ws.Range("M32").Value = -16959855
$ws.Range("N32").Value = -18654.916

$ws.Range("H61").Value = 7016.6
$ws.Range("I61").Value = 6277.75
$ws.Range("J61").Value = 8124.875
$ws.Range("K61").Value = 6277.75
$ws.Range("L61").Value = 8124.875
$ws.Range("M61").Value = -6065.75
$ws.Range("N61").Value = -8548.875

$ws.Range("H136").Value = 7016.6
$ws.Range("I136").Value = 6277.75
$ws.Range("J136").Value = 8124.875
$ws.Range("K136").Value = 18833.25
$ws.Range("L136").Value = 24374.625
$ws.Range("M136").Value = -16283.25
$ws.Range("N136").Value = -29474.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5716204.5
$ws.Range("I94").Value = 1446.4193
$ws.Range("K94").Value = 1446.4193
$ws.Range("M94").Value = -995.4193

$ws.Range("H134").Value = 2900.6858
$ws.Range("I134").Value = 1444.64
$ws.Range("K134").Value = 4333.92
$ws.Range("M134").Value = -1798.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1369.6666
$ws.Range("I2").Value = 1052
$ws.Range("K2").Value = 1052
$ws.Range("M2").Value = -939

$ws.Range("H12").Value = 584.3333
$ws.Range("I12").Value = 695.5
$ws.Range("J12").Value = 528.75
$ws.Range("K12").Value = 695.5
$ws.Range("L12").Value = 528.75
$ws.Range("M12").Value = -525.5
$ws.Range("N12").Value = -868.75

$ws.Range("H15").Value = 569.6667
$ws.Range("J15").Value = 569.6667
$ws.Range("L15").Value = 569.6667
$ws.Range("N15").Value = -909.6667

$ws.Range("H21").Value = 3666.6667
$ws.Range("I21").Value = 3000
$ws.Range("K21").Value = 3000
$ws.Range("M21").Value = -2765

$ws.Range("H48").Value = 42739
$ws.Range("J48").Value = 42739
$ws.Range("L48").Value = 42739
$ws.Range("N48").Value = -43691

$ws.Range("H58").Value = 1550
$ws.Range("I58").Value = 1005.2353
$ws.Range("J58").Value = 2873
$ws.Range("K58").Value = 1005.2353
$ws.Range("L58").Value = 2873
$ws.Range("M58").Value = -802.2353000000001
$ws.Range("N58").Value = -3279

$ws.Range("H60").Value = 5999.5
$ws.Range("I60").Value = 5999.5
$ws.Range("K60").Value = 5999.5
$ws.Range("M60").Value = -5488.5

$ws.Range("H122").Value = 11059757
$ws.Range("J122").Value = 7043.5
$ws.Range("L122").Value = 21130.5
$ws.Range("N122").Value = -26030.5

$ws.Range("H132").Value = 3805.3462
$ws.Range("I132").Value = 3655.3914
$ws.Range("K132").Value = 10966.1742
$ws.Range("M132").Value = -8436.174199999999

$ws.Range("H135").Value = 69322
$ws.Range("J135").Value = 69322
$ws.Range("L135").Value = 69322
$ws.Range("N135").Value = -79462

$ws.Range("H136").Value = 1550
$ws.Range("I136").Value = 1005.2353
$ws.Range("J136").Value = 2873
$ws.Range("K136").Value = 3015.7059
$ws.Range("L136").Value = 8619
$ws.Range("M136").Value = -465.7058999999999
$ws.Range("N136").Value = -13719

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 748.44446
$ws.Range("I8").Value = 748.44446
$ws.Range("K8").Value = 2245.33338
$ws.Range("M8").Value = -2106.33338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 14706858
$ws.Range("I97").Value = 776.2308
$ws.Range("J97").Value = 62501624
$ws.Range("K97").Value = 776.2308
$ws.Range("L97").Value = 62501624
$ws.Range("M97").Value = -280.2308
$ws.Range("N97").Value = -62502616

$ws.Range("H136").Value = 20880.666
$ws.Range("J136").Value = 20880.666
$ws.Range("L136").Value = 62641.99800000001
$ws.Range("N136").Value = -67741.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 61612.367
$ws.Range("I100").Value = 93686.664
$ws.Range("K100").Value = 93686.664
$ws.Range("M100").Value = -93145.664

$ws.Range("H132").Value = 4859.82
$ws.Range("I132").Value = 4653.1714
$ws.Range("K132").Value = 13959.5142
$ws.Range("M132").Value = -11429.5142

$ws.Range("H136").Value = 4927.423
$ws.Range("I136").Value = 2975.875
$ws.Range("J136").Value = 8049.9
$ws.Range("K136").Value = 8927.625
$ws.Range("L136").Value = 24149.7
$ws.Range("M136").Value = -6377.625
$ws.Range("N136").Value = -29249.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 42500
$ws.Range("I60").Value = 42500
$ws.Range("K60").Value = 42500
$ws.Range("M60").Value = -41678

$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

$ws.Range("H136").Value = 2472.4138
$ws.Range("I136").Value = 1828.6
$ws.Range("J136").Value = 6496.25
$ws.Range("K136").Value = 5485.799999999999
$ws.Range("L136").Value = 19488.75
$ws.Range("M136").Value = -2935.799999999999
$ws.Range("N136").Value = -24588.75
